$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells keep their existing text ("0.00"-style) representation
# rather than being auto-converted to numbers when we write values that
# look numeric.
$ws.Range("B3:B7").NumberFormat = "@"
$ws.Range("D3:D7").NumberFormat = "@"

# Row 3 (Glayds  Bundotich): 4.00 -> 1.00
$ws.Range("B3").Value = "1.00"
$ws.Range("D3").Value = "1.00"

# Row 4 (Jane Gichohi): 2.00 -> 1.00
$ws.Range("B4").Value = "1.00"
$ws.Range("D4").Value = "1.00"

# Row 5 (Mirriam Makau): 2.00 -> 1.00
$ws.Range("B5").Value = "1.00"
$ws.Range("D5").Value = "1.00"

# Row 6 (Victor Njogu): 2.00 -> 1.00
$ws.Range("B6").Value = "1.00"
$ws.Range("D6").Value = "1.00"

# Row 7 (KD Totals): 11.00 -> 5.00
$ws.Range("B7").Value = "5.00"
$ws.Range("D7").Value = "5.00"
